# Add a "SKU" column to the digital inventory header row.
# Inserting a new column B pushes the existing serial/name/location/inventory
# columns one place to the right (B->C, C->D, D->E, E->F) and picks up the
# formatting of the neighboring header cells automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(2).Insert()
$ws.Range("B1").Value = "SKU"

# Match the author's recorded selection on the newly added header cell.
$ws.Range("B1").Select()
